# Update the Rspo3-Lrp6 LR-pair sheet with the revised NATMI output
# (recomputed per Dr Hou's advice: added "sCs" as an additional sending
# cluster, doubling the number of data rows from 4 to 8).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,("FAPs", "Rspo3", "Lrp6", "ECs", 3, 1, 4.83236, 14.49708, 0.975350813525687, 0.975350813525687, 3, 1, 14.08030833333333, 42.240925, 0.1556474304006564, 0.1556474304006564, 68.04111877766667, 612.370068999, 0.151810847864463, 0.151810847864463)
    ,("FAPs", "Rspo3", "Lrp6", "FAPs", 3, 1, 4.83236, 14.49708, 0.975350813525687, 0.975350813525687, 3, 1, 43.15322233333333, 129.459667, 0.4770270657916382, 0.4770270657916382, 208.5319054747067, 1876.78714927236, 0.4652687366936457, 0.4652687366936458)
    ,("FAPs", "Rspo3", "Lrp6", "M2", 3, 1, 4.83236, 14.49708, 0.975350813525687, 0.975350813525687, 3, 1, 21.909414, 65.72824199999999, 0.2421924229181179, 0.242192422918118, 105.87417583704, 952.86758253336, 0.2362225767229436, 0.2362225767229436)
    ,("FAPs", "Rspo3", "Lrp6", "sCs", 3, 1, 4.83236, 14.49708, 0.975350813525687, 0.975350813525687, 3, 1, 11.31989366666667, 33.959681, 0.1251330808895874, 0.1251330808895874, 54.70180135905334, 492.3162122314801, 0.1220486522446347, 0.1220486522446347)
    ,("sCs", "Rspo3", "Lrp6", "ECs", 2, 0.6666666666666666, 0.122124, 0.366372, 0.02464918647431296, 0.02464918647431296, 3, 1, 14.08030833333333, 42.240925, 0.1556474304006564, 0.1556474304006564, 1.7195435749, 15.4758921741, 0.003836582536193428, 0.003836582536193429)
    ,("sCs", "Rspo3", "Lrp6", "FAPs", 2, 0.6666666666666666, 0.122124, 0.366372, 0.02464918647431296, 0.02464918647431296, 3, 1, 43.15322233333333, 129.459667, 0.4770270657916382, 0.4770270657916382, 5.270044124236, 47.430397118124, 0.01175832909799245, 0.01175832909799245)
    ,("sCs", "Rspo3", "Lrp6", "M2", 2, 0.6666666666666666, 0.122124, 0.366372, 0.02464918647431296, 0.02464918647431296, 3, 1, 21.909414, 65.72824199999999, 0.2421924229181179, 0.242192422918118, 2.675665275336, 24.080987478024, 0.005969846195174357, 0.005969846195174358)
    ,("sCs", "Rspo3", "Lrp6", "sCs", 2, 0.6666666666666666, 0.122124, 0.366372, 0.02464918647431296, 0.02464918647431296, 3, 1, 11.31989366666667, 33.959681, 0.1251330808895874, 0.1251330808895874, 1.382430694148, 12.441876247332, 0.003084428644952727, 0.003084428644952728)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowVals[$j]
    }
}
